$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old emoji "statut" codes to their new plain-text replacements
# (fix for Excel mangling/dropping emoji glyphs in the publipostage source).
$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📗" = "✅"
    "📙" = "+3"
}

# Values that look numeric ("-3", "+3") get silently coerced to actual
# numbers by the normal .Value setter, so those writes are done through a
# temporary Text number format and then restored to the default "Normal"
# style so the cell's formatting is left exactly as it started.
$numericLooking = @{ "-3" = $true; "+3" = $true }

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        if ($numericLooking.ContainsKey($newVal)) {
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newVal
        }
    }
}
